# Atualização de bases das ligas, do dia: 03-03-2024 às 00:35
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: push current row 103 (match id 101, Hyderabad FC vs Northeast United,
#     not yet played) down to a new row 104, keeping all its data/format intact,
#     only its running "id" (column A) changes from 101 to 102.
$ws.Range("A103:AC103").Copy($ws.Range("A104:AC104"))
$ws.Cells.Item(104, 1).Value = 102
# the source row had no values in H:J (score/result) or AB:AC (AH over/under PL) -
# the copy materialised empty placeholder cells there, so drop them again
$ws.Range("H104:J104").ClearContents()
$ws.Range("AB104:AC104").ClearContents()

# --- Step 2: row 103 becomes a brand-new fixture (Chennaiyin FC vs Odisha FC)
$ws.Cells.Item(103, 2).Value = 7749768
$ws.Cells.Item(103, 5).Value = 45354.45833333334
$ws.Cells.Item(103, 6).Value = "Chennaiyin FC"
$ws.Cells.Item(103, 7).Value = "Odisha FC"
$ws.Cells.Item(103, 11).Value = 2.6
$ws.Cells.Item(103, 12).Value = 3.5
$ws.Cells.Item(103, 13).Value = 2.25
$ws.Cells.Item(103, 14).Value = 2.9
$ws.Cells.Item(103, 15).Value = 3.5
$ws.Cells.Item(103, 16).Value = 2.05
$ws.Cells.Item(103, 17).Value = 0.25
$ws.Cells.Item(103, 18).Value = 1.9
$ws.Cells.Item(103, 19).Value = 1.9
$ws.Cells.Item(103, 20).Value = 2.75
$ws.Cells.Item(103, 21).Value = 2
$ws.Cells.Item(103, 22).Value = 1.8
$ws.Cells.Item(103, 23).Value = 0
$ws.Cells.Item(103, 24).Value = 0
$ws.Cells.Item(103, 25).Value = 0
$ws.Cells.Item(103, 26).Value = 0
$ws.Cells.Item(103, 27).Value = 0

# --- Step 3: row 101 (match id 99, Punjab FC vs Mumbai City FC) now has a final
#     score and closing odds
$ws.Cells.Item(101, 8).Value = 2
$ws.Cells.Item(101, 9).Value = 3
$ws.Cells.Item(101, 10).Value = "A"
$ws.Cells.Item(101, 14).Value = 3.4
$ws.Cells.Item(101, 15).Value = 3.4
$ws.Cells.Item(101, 16).Value = 2
$ws.Cells.Item(101, 18).Value = 1.8
$ws.Cells.Item(101, 19).Value = 2.05
$ws.Cells.Item(101, 20).Value = 2.5
$ws.Cells.Item(101, 21).Value = 1.825
$ws.Cells.Item(101, 22).Value = 2.025
$ws.Cells.Item(101, 23).Value = -1
$ws.Cells.Item(101, 24).Value = -1
$ws.Cells.Item(101, 25).Value = 1
$ws.Cells.Item(101, 26).Value = -1
$ws.Cells.Item(101, 27).Value = 1.05
$ws.Cells.Item(101, 28).Value = 0.825
$ws.Cells.Item(101, 29).Value = -1

# --- Step 4: row 102 (match id 100, Bengaluru vs Kerala Blasters) now has a final
#     score and closing odds
$ws.Cells.Item(102, 8).Value = 1
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = "H"
$ws.Cells.Item(102, 14).Value = 2.25
$ws.Cells.Item(102, 15).Value = 3.1
$ws.Cells.Item(102, 16).Value = 3.1
$ws.Cells.Item(102, 18).Value = 1.975
$ws.Cells.Item(102, 19).Value = 1.875
$ws.Cells.Item(102, 20).Value = 2.25
$ws.Cells.Item(102, 21).Value = 1.925
$ws.Cells.Item(102, 22).Value = 1.925
$ws.Cells.Item(102, 23).Value = 1.25
$ws.Cells.Item(102, 24).Value = -1
$ws.Cells.Item(102, 25).Value = -1
$ws.Cells.Item(102, 26).Value = 0.9750000000000001
$ws.Cells.Item(102, 27).Value = -1
$ws.Cells.Item(102, 28).Value = -1
$ws.Cells.Item(102, 29).Value = 0.925
